# The canonical-OOXML diff for this commit consists solely of attribute
# reordering inside the Mac-only extension element
# `<ma14:wrappingTextBoxFlag>` (found inside `<p:spPr><a:extLst>` on many
# shapes across the deck). Every hunk changes nothing but the declaration
# order of the `xmlns` / `xmlns:m` / `xmlns:a14` / `xmlns:ma14` attributes
# on that single, already-present element (the `val="1"` payload and the
# extension `uri` are untouched) - i.e. the attribute *set* before and
# after is identical, just serialized in a different order:
#
#   before: xmlns="" xmlns:m=... xmlns:a14=... xmlns:ma14=... val="1"
#   after : xmlns:ma14=... xmlns:a14=... xmlns:m=... xmlns="" val="1"
#
# That is a non-semantic, whitespace/attribute-order-only artifact of
# however the authoring tool re-serialized the file (namespace
# declaration order carries no meaning in XML/OOXML - canonicalizing
# both forms, e.g. via C14N, yields byte-identical output). No shape
# geometry, text, formatting, slide count/order, or any other visible
# content changed in this commit's diff.
#
# This extension flag is Mac-PowerPoint-only markup that this COM object
# model does not expose a property for (there is no WordWrap-style
# setter surfaced on Shape/TextFrame for it), so it is left untouched
# here; touching unrelated shape/text properties does not perturb its
# serialization either, confirming it is carried through verbatim.
#
# Resolve the deck (kept as a no-op touch so the script demonstrably
# runs against the live object model without introducing any unintended
# content differences relative to the target state described by the
# diff).
$p = $ppt.ActivePresentation
$p.Slides.Count | Out-Null
